$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Marking row (B11): number right marks multiplier
$ws.Range("B11").Value = 5

# Update the Total row (B12): total score
$ws.Range("B12").Value = 120

# Update the corr/total marks summary text (E12)
$ws.Range("E12").Value = "120/140"
